$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44204
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 7000
$ws.Cells.Item(2, 12).Value = 7500
$ws.Cells.Item(2, 13).Value = 7188
$ws.Cells.Item(2, 16).Value = 719

# Row 3
$ws.Cells.Item(3, 4).Value = 44194
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 8000
$ws.Cells.Item(3, 12).Value = 9000
$ws.Cells.Item(3, 13).Value = 8500
$ws.Cells.Item(3, 16).Value = 850

# Row 4
$ws.Cells.Item(4, 4).Value = 44428
$ws.Cells.Item(4, 10).Value = 50
$ws.Cells.Item(4, 11).Value = 7500
$ws.Cells.Item(4, 12).Value = 8000
$ws.Cells.Item(4, 13).Value = 7800
$ws.Cells.Item(4, 16).Value = 780

# Row 5
$ws.Cells.Item(5, 4).Value = 44362
$ws.Cells.Item(5, 10).Value = 50
$ws.Cells.Item(5, 11).Value = 6000
$ws.Cells.Item(5, 12).Value = 6500
$ws.Cells.Item(5, 13).Value = 6300
$ws.Cells.Item(5, 16).Value = 630

# Row 6
$ws.Cells.Item(6, 4).Value = 44421
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 7500
$ws.Cells.Item(6, 13).Value = 7250
$ws.Cells.Item(6, 16).Value = 725

# Row 7
$ws.Cells.Item(7, 4).Value = 44349
$ws.Cells.Item(7, 11).Value = 6000
$ws.Cells.Item(7, 12).Value = 6500
$ws.Cells.Item(7, 13).Value = 6250
$ws.Cells.Item(7, 16).Value = 625

# Row 8
$ws.Cells.Item(8, 4).Value = 44231
$ws.Cells.Item(8, 10).Value = 70
$ws.Cells.Item(8, 11).Value = 7500
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 7714
$ws.Cells.Item(8, 16).Value = 771

# Row 9
$ws.Cells.Item(9, 4).Value = 44435
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 7000
$ws.Cells.Item(9, 12).Value = 7500
$ws.Cells.Item(9, 13).Value = 7250
$ws.Cells.Item(9, 16).Value = 725

# Row 10
$ws.Cells.Item(10, 4).Value = 44313
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 6000
$ws.Cells.Item(10, 12).Value = 6500
$ws.Cells.Item(10, 13).Value = 6250
$ws.Cells.Item(10, 16).Value = 625

# Row 11
$ws.Cells.Item(11, 4).Value = 44188
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 8500
$ws.Cells.Item(11, 13).Value = 8250
$ws.Cells.Item(11, 16).Value = 825

# Row 12
$ws.Cells.Item(12, 4).Value = 44230
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = 9500
$ws.Cells.Item(12, 16).Value = 950

# Row 13
$ws.Cells.Item(13, 4).Value = 44355
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 6000
$ws.Cells.Item(13, 12).Value = 6500
$ws.Cells.Item(13, 13).Value = 6300
$ws.Cells.Item(13, 16).Value = 630

# Row 14
$ws.Cells.Item(14, 4).Value = 44238
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 8000
$ws.Cells.Item(14, 12).Value = 8500
$ws.Cells.Item(14, 13).Value = 8250
$ws.Cells.Item(14, 16).Value = 825

# Row 15
$ws.Cells.Item(15, 4).Value = 44299
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 16).Value = 850

# Row 16
$ws.Cells.Item(16, 4).Value = 44320
$ws.Cells.Item(16, 10).Value = 50
$ws.Cells.Item(16, 11).Value = 7000
$ws.Cells.Item(16, 12).Value = 7500
$ws.Cells.Item(16, 13).Value = 7200
$ws.Cells.Item(16, 16).Value = 720

# Row 17
$ws.Cells.Item(17, 4).Value = 44334
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 6500
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 6750
$ws.Cells.Item(17, 16).Value = 675

# Row 18
$ws.Cells.Item(18, 4).Value = 44342
$ws.Cells.Item(18, 10).Value = 50
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 6500
$ws.Cells.Item(18, 13).Value = 6300
$ws.Cells.Item(18, 16).Value = 630

# Row 19
$ws.Cells.Item(19, 4).Value = 44225
$ws.Cells.Item(19, 11).Value = 7500
$ws.Cells.Item(19, 12).Value = 8000
$ws.Cells.Item(19, 13).Value = 7750
$ws.Cells.Item(19, 16).Value = 775

# Row 20
$ws.Cells.Item(20, 4).Value = 44327

# Row 21
$ws.Cells.Item(21, 4).Value = 44399
$ws.Cells.Item(21, 10).Value = 60
$ws.Cells.Item(21, 11).Value = 9000
$ws.Cells.Item(21, 12).Value = 10000
$ws.Cells.Item(21, 13).Value = 9500
$ws.Cells.Item(21, 16).Value = 950

# Row 22
$ws.Cells.Item(22, 4).Value = 44383
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 13).Value = 7750
$ws.Cells.Item(22, 16).Value = 775

# Row 23
$ws.Cells.Item(23, 4).Value = 44405
$ws.Cells.Item(23, 10).Value = 80
$ws.Cells.Item(23, 11).Value = 7500
$ws.Cells.Item(23, 12).Value = 8000
$ws.Cells.Item(23, 13).Value = 7688
$ws.Cells.Item(23, 16).Value = 769

# Row 24
$ws.Cells.Item(24, 4).Value = 44246
$ws.Cells.Item(24, 10).Value = 60
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 9500
$ws.Cells.Item(24, 16).Value = 950

# Row 25
$ws.Cells.Item(25, 4).Value = 44336
$ws.Cells.Item(25, 10).Value = 60
$ws.Cells.Item(25, 11).Value = 6000
$ws.Cells.Item(25, 12).Value = 6500
$ws.Cells.Item(25, 13).Value = 6250
$ws.Cells.Item(25, 16).Value = 625

# Row 26
$ws.Cells.Item(26, 4).Value = 44328
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 7500
$ws.Cells.Item(26, 13).Value = 7250
$ws.Cells.Item(26, 16).Value = 725

# Row 27
$ws.Cells.Item(27, 4).Value = 44453
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 6500
$ws.Cells.Item(27, 12).Value = 7000
$ws.Cells.Item(27, 13).Value = 6750
$ws.Cells.Item(27, 16).Value = 675

# Row 28
$ws.Cells.Item(28, 4).Value = 44308
$ws.Cells.Item(28, 11).Value = 5000
$ws.Cells.Item(28, 12).Value = 5500
$ws.Cells.Item(28, 13).Value = 5250
$ws.Cells.Item(28, 16).Value = 525

# Row 29
$ws.Cells.Item(29, 4).Value = 44264
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 8500
$ws.Cells.Item(29, 13).Value = 8200
$ws.Cells.Item(29, 16).Value = 820

# Row 30
$ws.Cells.Item(30, 4).Value = 44166

# Row 31
$ws.Cells.Item(31, 4).Value = 44433
$ws.Cells.Item(31, 10).Value = 100
$ws.Cells.Item(31, 11).Value = 7000
$ws.Cells.Item(31, 12).Value = 7500
$ws.Cells.Item(31, 13).Value = 7250
$ws.Cells.Item(31, 16).Value = 725

# Row 33
$ws.Cells.Item(33, 4).Value = 44279
$ws.Cells.Item(33, 10).Value = 60
$ws.Cells.Item(33, 11).Value = 7500
$ws.Cells.Item(33, 12).Value = 8000
$ws.Cells.Item(33, 13).Value = 7750
$ws.Cells.Item(33, 16).Value = 775

# Row 35
$ws.Cells.Item(35, 4).Value = 44209
$ws.Cells.Item(35, 10).Value = 80
$ws.Cells.Item(35, 11).Value = 7500
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 7688
$ws.Cells.Item(35, 16).Value = 769

# Row 36
$ws.Cells.Item(36, 4).Value = 44160
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 9000
$ws.Cells.Item(36, 12).Value = 9500
$ws.Cells.Item(36, 13).Value = 9250
$ws.Cells.Item(36, 16).Value = 925

# Row 37
$ws.Cells.Item(37, 4).Value = 44351
$ws.Cells.Item(37, 11).Value = 6000
$ws.Cells.Item(37, 12).Value = 6500
$ws.Cells.Item(37, 13).Value = 6300
$ws.Cells.Item(37, 16).Value = 630

# Row 38
$ws.Cells.Item(38, 4).Value = 44365
$ws.Cells.Item(38, 10).Value = 50
$ws.Cells.Item(38, 11).Value = 6000
$ws.Cells.Item(38, 12).Value = 6500
$ws.Cells.Item(38, 13).Value = 6200
$ws.Cells.Item(38, 16).Value = 620

# Row 39
$ws.Cells.Item(39, 4).Value = 44357
$ws.Cells.Item(39, 10).Value = 50
$ws.Cells.Item(39, 13).Value = 6200
$ws.Cells.Item(39, 16).Value = 620

# Row 40
$ws.Cells.Item(40, 4).Value = 44376
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 6000
$ws.Cells.Item(40, 12).Value = 6500
$ws.Cells.Item(40, 13).Value = 6250
$ws.Cells.Item(40, 16).Value = 625

# Row 41
$ws.Cells.Item(41, 4).Value = 44292
$ws.Cells.Item(41, 10).Value = 50
$ws.Cells.Item(41, 11).Value = 10000
$ws.Cells.Item(41, 12).Value = 11000
$ws.Cells.Item(41, 13).Value = 10600
$ws.Cells.Item(41, 16).Value = 1060

# Row 42
$ws.Cells.Item(42, 4).Value = 44425
$ws.Cells.Item(42, 11).Value = 6500
$ws.Cells.Item(42, 12).Value = 7000
$ws.Cells.Item(42, 13).Value = 6750
$ws.Cells.Item(42, 16).Value = 675

